# update contract details ui
#
# The "Адрес, телефон" / "{deadman_address}" table row in the contract
# header actually held the client's address placeholder in the wrong
# cell. This fixes it so that:
#   - the label cell keeps its caption and now also carries the
#     deadman's address placeholder inline, and
#   - the value cell (which used to wrongly repeat {deadman_address})
#     now correctly holds the client's phone placeholder.
#
# Order matters: we retarget the existing lone {deadman_address} run to
# {client_phone} FIRST, while it is still unique in the document, and
# only THEN introduce the new {deadman_address} token into the label
# cell's text. Doing it the other way around would leave two
# {deadman_address} runs in the doc at the same time and the second
# Find/Replace could clobber the wrong one.

$d = $word.ActiveDocument

# 1) The old (misplaced) {deadman_address} value becomes {client_phone}.
$d.Content.Find.Execute(
    "{deadman_address}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{client_phone}", 2) | Out-Null

# 2) The "Адрес, телефон" label run gains the real {deadman_address}
#    placeholder in place of its long run of trailing spaces.
$d.Content.Find.Execute(
    "Адрес, телефон                                       ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Адрес, телефон  {deadman_address}", 2) | Out-Null
